# Trade #6 (overall trade #67, MarketMaking strategy trade #35) closed at
# 2026-02-17 20:48:11 - unknown UNKNOWN +0.000%
#
# This script:
#   1. Updates the aggregate metrics on the "Summary" sheet.
#   2. Updates the MarketMaking row on the "Strategy Status" sheet.
#   3. Marks the open MarketMaking trade (row 35 on "All Trades" /
#      row 2 on "MarketMaking") as CLOSED with its exit data.
#   4. Appends a brand new OPEN MarketMaking trade row to both the
#      "All Trades" sheet (row 68) and the "MarketMaking" sheet (row 35).

$wb = $excel.ActiveWorkbook

# Helper so text that looks like a date/time ("2026-02-17") is written as
# plain text instead of being auto-converted into a date serial number by
# Excel's auto-detection, while leaving the cell's style untouched
# afterwards (matches the workbook's original unstyled inline strings).
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1400.3             # Current Capital
$summary.Range("B4").Value = 0.09               # Total P&L $
$summary.Range("B5").Value = 0.05               # Total P&L %
$summary.Range("B6").Value = 34                 # Total Trades
$summary.Range("B8").Value = 15                 # Losing Trades
$summary.Range("B9").Value = 41.18              # Win Rate %

# ---------------------------------------------------------------------
# 2. Strategy Status sheet (row 5 = MarketMaking)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 100.3               # Capital
$status.Range("D5").Value = 1                   # Trades
$status.Range("E5").Value = -0.02               # P&L $
$status.Range("F5").Value = 0.3                 # P&L %

# ---------------------------------------------------------------------
# 3. All Trades sheet - close the open trade (row 35)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G35").Value = 0.090909        # Exit Price
$allTrades.Range("H35").Value = "CLOSED"        # Status
$allTrades.Range("I35").Value = -17.3554        # P&L %
$allTrades.Range("J35").Value = -0.02           # P&L $
$allTrades.Range("K35").Value = 100.3           # Capital After
$allTrades.Range("L35").Value = "early_exit"    # Exit Reason
$allTrades.Range("M35").Value = 0.13            # Duration (min)

# ---------------------------------------------------------------------
# 4. All Trades sheet - append the new open trade (row 68)
# ---------------------------------------------------------------------
$allTrades.Cells.Item(68, 1).Value = 67                 # Trade #
Set-TextValue $allTrades.Cells.Item(68, 2) "2026-02-17" # Date
Set-TextValue $allTrades.Cells.Item(68, 3) "20:48:04"   # Time
$allTrades.Cells.Item(68, 4).Value = "MarketMaking"     # Strategy
$allTrades.Cells.Item(68, 5).Value = "UP"                # Side
$allTrades.Cells.Item(68, 6).Value = 0.11                # Entry Price
$allTrades.Cells.Item(68, 8).Value = "OPEN"               # Status
$allTrades.Cells.Item(68, 9).Value = 0                    # P&L %
$allTrades.Cells.Item(68, 10).Value = 0                   # P&L $
$allTrades.Cells.Item(68, 11).Value = 100.32               # Capital After
$allTrades.Cells.Item(68, 13).Value = 0                    # Duration (min)
$allTrades.Cells.Item(68, 14).Value = 0                    # Entry Slippage (bps)
$allTrades.Cells.Item(68, 15).Value = 0                    # Exit Slippage (bps)
$allTrades.Cells.Item(68, 16).Value = 0.6                  # Confidence
$allTrades.Cells.Item(68, 17).Value = "Normal spread capture: 19600 bps" # Entry Reason

# ---------------------------------------------------------------------
# 5. MarketMaking sheet - close the open trade (row 2)
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("G2").Value = 0.090909                # Exit Price
$mm.Range("H2").Value = "CLOSED"                # Status
$mm.Range("I2").Value = -17.3554                # P&L %
$mm.Range("J2").Value = -0.02                   # P&L $
$mm.Range("K2").Value = 100.3                   # Capital After
$mm.Range("P2").Value = "early_exit"            # Exit Reason
$mm.Range("Q2").Value = 0.13                    # Duration (min)

# ---------------------------------------------------------------------
# 6. MarketMaking sheet - append the new open trade (row 35)
# ---------------------------------------------------------------------
$mm.Cells.Item(35, 1).Value = 67                         # Trade #
Set-TextValue $mm.Cells.Item(35, 2) "2026-02-17"         # Date
Set-TextValue $mm.Cells.Item(35, 3) "20:48:04"           # Time
$mm.Cells.Item(35, 4).Value = "MarketMaking"              # Strategy
$mm.Cells.Item(35, 5).Value = "UP"                         # Side
$mm.Cells.Item(35, 6).Value = 0.11                         # Entry Price
$mm.Cells.Item(35, 8).Value = "OPEN"                        # Status
$mm.Cells.Item(35, 9).Value = 0                             # P&L %
$mm.Cells.Item(35, 10).Value = 0                            # P&L $
$mm.Cells.Item(35, 11).Value = 100.32                       # Capital After
$mm.Cells.Item(35, 12).Value = 0                            # Entry Slippage (bps)
$mm.Cells.Item(35, 13).Value = 0                            # Exit Slippage (bps)
$mm.Cells.Item(35, 14).Value = 0.6                          # Confidence
$mm.Cells.Item(35, 15).Value = "Normal spread capture: 19600 bps" # Entry Reason
$mm.Cells.Item(35, 17).Value = 0                            # Duration (min)
